$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.124.02"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Value = "2.634.48"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.95"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.66"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.335"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("D13").Value = "3.109.04"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("D14").Value = "59.123.52"
$ws.Range("E14").Value = "  -0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.64"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.53%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000136"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.609.56"
$ws.Range("E17").Value = "  -2.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "342.90"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.41"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.52"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.34"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.42"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.413"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("D26").Value = "2.760.84"
$ws.Range("E26").Value = "  -0.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.14"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.53%  "

$ws.Range("D29").Value = "0.0₃0790"
$ws.Range("E29").Value = "  -1.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.60"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.96"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.91"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.12"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.16"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.850"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.24"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.839"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.72%  "

$ws.Range("E40").Value = "  -2.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.61"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0976"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.596"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.72%  "

$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.71"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.70%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "268.17"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.10"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0531"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.11%  "

$ws.Range("D49").Value = "2.033.71"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0228"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.72"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.19%  "
